# "type split calibration - i think its good enough now"
#
# Re-calibrates the dwelling-type split for years 1901-2050 (rows 303-452):
#   - columns B:F (pre-1955 .. 1991-2000 construction-era buckets) are scaled
#     down by a factor of 0.8
#   - columns G:H (2001-2010, post-2010 buckets) are scaled up by a factor
#     of 1.21
#
# Zero-valued cells remain zero either way, so the whole B:H block for the
# row range can be rescaled uniformly without needing to special-case the
# "not yet built up" columns for the earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 303
$lastRow  = 452

$rng = $ws.Range("B$($firstRow):H$($lastRow)")
$vals = $rng.Value2

$rowCount = $vals.GetLength(0)
$colCount = $vals.GetLength(1)

for ($i = 1; $i -le $rowCount; $i++) {
    for ($j = 1; $j -le $colCount; $j++) {
        if ($j -le 5) {
            # columns B,C,D,E,F
            $factor = 0.8
        } else {
            # columns G,H
            $factor = 1.21
        }
        $vals[$i, $j] = $vals[$i, $j] * $factor
    }
}

$rng.Value2 = $vals
